$wb = $excel.ActiveWorkbook

# --- Sheets (tab order: 1=Overview, 2=zh-cn, 3=de-de) ---
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# --- Text change: "Ready for handoff" -> "In Translation" everywhere it appears ---
# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows (2, 3)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width changes: narrow the "Ready for handoff"/status columns ---
# Overview: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de: column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
